# Applies the "Added quantities to products." commit.
#
# On the "Products" sheet:
#  - fills in the new Quantity column (E2:E49)
#  - a cosmetic re-fit of column D (Price) width
#  - the sheet view left at the zoom/scroll/selection the author
#    had when they saved (zoom 85%, scrolled to row 31, F44 selected)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Products")
$ws.Activate()

# Quantity values for column E, row 2 through 49, in sheet order.
$quantities = @(
    24, 15, 8, 30, 26, 23, 33, 6, 18, 32, 21, 31, 7, 12, 18, 12, 29, 15,
    2, 0, 3, 18, 22, 1, 3, 2, 30, 3, 4, 5, 21, 9, 23, 21, 9, 7, 2, 2, 29,
    18, 0, 28, 12, 14, 16, 32, 8, 18
)

$firstRow = 2
for ($i = 0; $i -lt $quantities.Count; $i++) {
    $ws.Cells.Item($firstRow + $i, 5).Value = $quantities[$i]
}

# Column D (Price) width: author nudged it slightly wider.
$ws.Columns.Item(4).ColumnWidth = 7.833333333333333

# Sheet view state captured at save time: zoom, scroll position, selection.
$excel.ActiveWindow.Zoom = 85
$excel.ActiveWindow.ScrollRow = 31
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("F44").Select()

